$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3138
$ws1.Range("F4").Value = 1960
$ws1.Range("F5").Value = 256
$ws1.Range("F6").Value = 78
$ws1.Range("F7").Value = 2625
$ws1.Range("F8").Value = 596
$ws1.Range("F14").Value = 9882
$ws1.Range("F16").Value = 26
$ws1.Range("F18").Value = 7804
$ws1.Range("F19").Value = 12397
$ws1.Range("F20").Value = 127
$ws1.Range("F24").Value = 586
$ws1.Range("F26").Value = 251
$ws1.Range("F27").Value = 227
$ws1.Range("F28").Value = 4188
$ws1.Range("F29").Value = 1339
$ws1.Range("F32").Value = 72
$ws1.Range("F33").Value = 4598
$ws1.Range("F34").Value = 1283
$ws1.Range("F35").Value = 60
$ws1.Range("F37").Value = 69
$ws1.Range("F38").Value = 611

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 646

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 646
$ws4.Range("F4").Value = 3138
$ws4.Range("F6").Value = 1960
$ws4.Range("F8").Value = 256
$ws4.Range("F9").Value = 2625
$ws4.Range("F11").Value = 596
$ws4.Range("F17").Value = 9882
$ws4.Range("F19").Value = 26
$ws4.Range("F21").Value = 7805
$ws4.Range("F22").Value = 12397
$ws4.Range("F24").Value = 127
$ws4.Range("F27").Value = 586
$ws4.Range("F32").Value = 251
$ws4.Range("F33").Value = 227
$ws4.Range("F36").Value = 72
$ws4.Range("F37").Value = 4598
$ws4.Range("F45").Value = 611
